$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.846.34'
$ws.Range("E2").Value = '  +1.02%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.707.50'
$ws.Range("E3").Value = '  +0.85%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.40%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.80'
$ws.Range("E5").Value = '  +0.84%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("E6").Value = '  +0.20%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3969'
$ws.Range("E7").Value = '  +0.61%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4105'
$ws.Range("E8").Value = '  +2.28%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.513'
$ws.Range("E9").Value = '  -0.48%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.006'
$ws.Range("E10").Value = '  +0.46%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.57'
$ws.Range("E11").Value = '  +1.91%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08981'
$ws.Range("E12").Value = '  +2.59%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.698'
$ws.Range("E13").Value = '  +6.72%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '24.39'
$ws.Range("E14").Value = '  +4.96%  '

$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.224'
$ws.Range("E15").Value = '  +0.38%  '

$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001372'
$ws.Range("E16").Value = '  +4.58%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.719.04'
$ws.Range("E17").Value = '  +1.26%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '100.32'
$ws.Range("E18").Value = '  +0.53%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07166'
$ws.Range("E19").Value = '  +1.39%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.554'
$ws.Range("E20").Value = '  +6.87%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '20.11'
$ws.Range("E21").Value = '  +2.35%  '

$ws.Range("E22").Value = '  +0.85%  '

$ws.Range("E23").Value = '  +2.71%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.848.37'
$ws.Range("E24").Value = '  +1.06%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.091'
$ws.Range("E25").Value = '  -1.16%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.345'
$ws.Range("E26").Value = '  +0.35%  '

$ws.Range("E27").Value = '  +1.37%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.348'
$ws.Range("E28").Value = '  +25.21%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '166.90'
$ws.Range("E29").Value = '  +3.17%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '139.51'

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.235'
$ws.Range("E31").Value = '  +0.76%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.901'
$ws.Range("E32").Value = '  +9.99%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09129'
$ws.Range("E33").Value = '  +6.14%  '

$ws.Range("E34").Value = '  +0.95%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.03052'
$ws.Range("E35").Value = '  +11.70%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.2823'
$ws.Range("E36").Value = '  +3.34%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '11.14'

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.966'
$ws.Range("E38").Value = '  +2.14%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '14.61'
$ws.Range("E39").Value = '  +1.49%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.09346'
$ws.Range("E40").Value = '  +2.53%  '

$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.7913'
$ws.Range("E41").Value = '  +3.75%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.484'
$ws.Range("E42").Value = '  -0.16%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '16.90'
$ws.Range("E43").Value = '  +7.97%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.654'
$ws.Range("E44").Value = '  +2.38%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.7338'
$ws.Range("E45").Value = '  +2.59%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.279'
$ws.Range("E46").Value = '  +1.40%  '

$ws.Range("B47").Value = 'Frax'
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.005'
$ws.Range("E47").Value = '  +0.54%  '

$ws.Range("B48").Value = 'Flow'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.352'
$ws.Range("E48").Value = '  +2.43%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.19'
$ws.Range("E49").Value = '  +0.26%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '94.79'
$ws.Range("E50").Value = '  +5.61%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.08078'
$ws.Range("E51").Value = '  +1.19%  '
